$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C2").Value = 26
$ws.Range("C3").Value = 27
$ws.Range("C4").Value = 18
$ws.Range("C5").Value = 28
$ws.Range("C6").Value = 21
$ws.Range("C7").Value = 19
$ws.Range("C8").Value = 23
$ws.Range("C9").Value = 31
$ws.Range("C10").Value = 28
$ws.Range("C11").Value = 28
$ws.Range("C12").Value = 36
$ws.Range("C13").Value = 29
$ws.Range("C14").Value = 18
$ws.Range("C15").Value = 24
$ws.Range("C16").Value = 28
$ws.Range("C17").Value = 19

# Update text values in column B
$ws.Range("B5").Value = "<fox>"
$ws.Range("B10").Value = "<see>"
$ws.Range("B13").Value = "<said>"
$ws.Range("B16").Value = "<he>"
$ws.Range("B17").Value = "<in>"
$ws.Range("B18").Value = "<he>"
